$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "Marking" row (row 11): negative marking was -1, should be -2,
# and the marks-per-right-answer was 5, should be 4.
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Correct the "Total" row (row 12): total obtained marks and max marks.
$ws.Range("B12").Value = 56
$ws.Range("E12").Value = "56 / 112"
